$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append rows 18-24 to the translation table, continuing the same layout as
# the existing rows. Formats are copied from whichever earlier row uses the
# matching style group (A6:E6 -> styles 4/5, A7:E7 -> styles 10/11) so the new
# cells reuse the existing style indices instead of creating new ones. Row 24
# has no filename (column A), mirroring row 14, so only B:E are touched there.

$ws.Range("A7:E7").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("A18").Value = 'SCRIPT/T01P02A/um2501.ssb'
$ws.Rows("18:18").RowHeight = 43.2

$ws.Range("A6:E6").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A19").Value = 'SCRIPT/D73P23A/us3103.ssb'
$ws.Range("B19").Value = 225
$ws.Range("C19").Value = ' I was thinking that [CS:I]Sky Gift[CR] was\nsome kind of a joke, because it was completely\nempty when I opened it.'
$ws.Range("D19").Value = ' Сначала я подумал, что [CS:I]Небесный\nПодарок[CR] это какой-то розыгрыш, потому что\nон был совершенно пуст.'
$ws.Range("E19").Value = ' Òîàœàìà ÿ ðïäôíàì, œóï [CS:I]Îåáåòîúê\nÐïäàñïë[CR] üóï ëàëïê-óï ñïèúãñúš, ðïóïíô œóï\nïî áúì òïâåñšåîîï ðôòó.'
$ws.Rows("19:19").RowHeight = 43.2

$ws.Range("A7:E7").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("B20").Value = 228
$ws.Range("C20").Value = ' But I was wrong.'
$ws.Range("D20").Value = ' Но я просто не понял сути.'
$ws.Range("E20").Value = ' Îï ÿ ðñïòóï îå ðïîÿì òôóé.'

$ws.Range("A6:E6").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A21").Value = 'SCRIPT/D73P26A/us3106.ssb'
$ws.Range("B21").Value = 212
$ws.Range("C21").Value = ' Looks like there will be snow\nfrom here on out.'
$ws.Range("D21").Value = ' Похоже, дальше будет много\nснега.'
$ws.Range("E21").Value = ' Ðïöïçå, äàìûšå áôäåó íîïãï\nòîåãà.'
$ws.Rows("21:21").RowHeight = 43.2

$ws.Range("A7:E7").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("B22").Value = 199
$ws.Range("C22").Value = ' It looks like the going will be\ntough from here on.'
$ws.Range("D22").Value = ' Похоже, дальше подъём будет\nтруден.'
$ws.Range("E22").Value = ' Ðïöïçå, äàìûšå ðïäùæí áôäåó\nóñôäåî.'
$ws.Rows("22:22").RowHeight = 21.6

$ws.Range("A6:E6").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Range("A23").Value = 'SCRIPT/P01P04A/us3114.ssb'
$ws.Range("B23").Value = 183
$ws.Range("C23").Value = ' Did you climb all the way to the\ntop of [CS:P]Sky Peak[CR]?'
$ws.Range("D23").Value = ' Вы сумели добраться до вершины\n[CS:P]Небесного Пика[CR]?'
$ws.Range("E23").Value = ' Âú òôíåìé äïáñàóûòÿ äï âåñšéîú\n[CS:P]Îåáåòîïãï Ðéëà[CR]?'
$ws.Rows("23:23").RowHeight = 43.2

$ws.Range("B6:E6").Copy()
$ws.Range("B24:E24").PasteSpecial(-4122)
$ws.Range("B24").Value = 186
$ws.Range("C24").Value = ' We got tired right around the\nsnowy area, so we decided to take a break.'
$ws.Range("D24").Value = ' Мы попытались пройти через\nзаснеженную местность, но устали и решили\nсделать перерыв.'
$ws.Range("E24").Value = ' Íú ðïðúóàìéòû ðñïêóé œåñåè\nèàòîåçåîîôý íåòóîïòóû, îï ôòóàìé é ñåšéìé\nòäåìàóû ðåñåñúâ.'
$ws.Rows("24:24").RowHeight = 31.8
